$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely - shifts B:F left to A:E, carrying values/styles with them
$ws.Columns.Item(1).Delete()
